$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'67.405.74"
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = "'3.310.32"
$ws.Range('E3').Value = '  +1.10%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'186.75"
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').Value = "'578.38"
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'0.605"
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('D10').Value = "'6.67"
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('D11').Value = "'0.408"
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('D12').Value = "'3.886.87"
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').Value = "'27.50"
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('D15').Value = "'67.658.73"
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = "'3.309.37"
$ws.Range('E17').Value = '  +2.12%  '
$ws.Range('D18').Value = "'445.06"
$ws.Range('E18').Value = '  +6.42%  '
$ws.Range('D19').Value = "'5.70"
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('D20').Value = "'13.60"
$ws.Range('E20').Value = '  +1.82%  '
$ws.Range('D21').Value = "'7.75"
$ws.Range('E21').Value = '  +2.54%  '
$ws.Range('D22').Value = "'73.97"
$ws.Range('E22').Value = '  +3.42%  '
$ws.Range('D23').Value = "'0.999"
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').Value = "'0.518"
$ws.Range('E24').Value = '  +1.88%  '
$ws.Range('D25').Value = "'3.455.35"
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('E26').Value = '  +1.26%  '
$ws.Range('E27').Value = '  +0.89%  '
$ws.Range('D29').Value = "'0.998"
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('D31').Value = "'22.94"
$ws.Range('E31').Value = '  +0.79%  '
$ws.Range('D32').Value = "'5.34"
$ws.Range('E32').Value = '  -2.43%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').Value = "'0.999"
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = "'1.25"
$ws.Range('E34').Value = '  -0.17%  '
$ws.Range('D35').Value = "'6.81"
$ws.Range('E35').Value = '  -1.27%  '
$ws.Range('E36').Value = '  +4.86%  '
$ws.Range('D37').Value = "'162.77"
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('D39').Value = "'27.21"
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').Value = "'0.791"
$ws.Range('E40').Value = '  -1.08%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').Value = "'2.757.20"
$ws.Range('E42').Value = '  +3.32%  '
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').Value = "'0.0676"
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('D45').Value = "'24.92"
$ws.Range('E45').Value = '  +1.59%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').Value = "'2.42"
$ws.Range('E46').Value = '  -0.56%  '
$ws.Range('D47').Value = "'40.21"
$ws.Range('E47').Value = '  -1.72%  '
$ws.Range('D48').Value = "'326.84"
$ws.Range('E48').Value = '  -3.19%  '
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').Value = "'0.992"
$ws.Range('E50').Value = '  +1.32%  '
$ws.Range('D51').Value = "'31.25"
$ws.Range('E51').Value = '  +1.43%  '
